# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-08-02 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2023-08-03 Thursday", 2) | Out-Null

# Update each arithmetic-expression cell in the table, in row-major
# document order. Using direct cell indexing (rather than text search)
# because several of the original expressions are not unique within
# the table (e.g. "88-26=" appears twice with different replacements).
$newValues = @(
    "99-55=",
    "83-3=",
    "76-48=",
    "57+34=",
    "80+9=",
    "84-15=",
    "9+67=",
    "62+9=",
    "86+11=",
    "67+15=",
    "50+12=",
    "91+4=",
    "97-81=",
    "59-48=",
    "65-25=",
    "47+23=",
    "38+7=",
    "88-82=",
    "89-10=",
    "91-45=",
    "67+0=",
    "51-37=",
    "29+47=",
    "31+52=",
    "18+39=",
    "22+16=",
    "77-65=",
    "29+0=",
    "27+69=",
    "63+7=",
    "80-18=",
    "45+43=",
    "81+4=",
    "82-3=",
    "27+53=",
    "62+14=",
    "14-13=",
    "69-37=",
    "33+14=",
    "42+56=",
    "40-39=",
    "37-18=",
    "29-6=",
    "49-30=",
    "73-55=",
    "58-45=",
    "28+61=",
    "15-3=",
    "67-26=",
    "2+46=",
    "11+84=",
    "45+53=",
    "23+22=",
    "94-71=",
    "56+12=",
    "24-2=",
    "6+18=",
    "11+0=",
    "76-64=",
    "33+66=",
    "42+38=",
    "72+12=",
    "73-25=",
    "50-17=",
    "51+40=",
    "79-52=",
    "21+62=",
    "43-21=",
    "45+44=",
    "95-77=",
    "36+17=",
    "87-87=",
    "95-56=",
    "28+48=",
    "62+23=",
    "93-34=",
    "60-31=",
    "7+1=",
    "98-67=",
    "8+61=",
    "49-36=",
    "48+2=",
    "41+23=",
    "75-72=",
    "27+63=",
    "93-8=",
    "11+88=",
    "94+2=",
    "76-23=",
    "50-37=",
    "34+9=",
    "21+46=",
    "27-9=",
    "80-67=",
    "51+3=",
    "42+38=",
    "26+61=",
    "15+27=",
    "78-47=",
    "81-37="
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i++
    }
}
